$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing cells
$ws.Range("C2").Value = "n/a"
$ws.Range("C3").Value = "n/a"
$ws.Range("A4").Value = 'Special Blue Pocket Torch,3"Dual Mini Plastic (20PC CT) '
$ws.Range("B7").Value = "https://mrawholesale.com/products/SPECIAL-BLUE-POCKET-TORCH-2-5-CLASSIC-RUBBER-LIGHTER-20CT-p572907865"
$ws.Range("C8").Value = "n/a"

# New rows data: A = product_name, B = product_link, C = brand, D = Flavors, E = url
$newRows = @(
    @("MAGNETIC TRAY WITH COVER SMALL", "https://mrawholesale.com/products/MAGNETIC-TRAY-WITH-COVER-SMALL-p572926723", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("MAGENTIC TRAY WITH 3D DESIGN WITH COVER SMALL", "https://mrawholesale.com/products/MAGENTIC-TRAY-WITH-3D-DESIGN-WITH-COVER-SMALL-p572927701", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("MAGNETIC TRAY MDEIUM SIZE ASSORTED", "https://mrawholesale.com/products/MAGNETIC-TRAY-MDEIUM-SIZE-ASSORTED-p572927736", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("MAGNETIC TRAY 3D DESIGN MIDIUM SIZE", "https://mrawholesale.com/products/MAGNETIC-TRAY-3D-DESIGN-MIDIUM-SIZE-p572927739", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("BLUNTLIFE 24CT LARGE INCENSE STICK", "https://mrawholesale.com/products/BLUNTLIFE-24CT-LARGE-INCENSE-STICK-p572929820", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("BLUNTLIFE 72CT SMALL INCENSE STICK", "https://mrawholesale.com/products/BLUNTLIFE-72CT-SMALL-INCENSE-STICK-p572927961", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("BLUNTLIFE 20CT SPRAY", "https://mrawholesale.com/products/BLUNTLIFE-20CT-SPRAY-p572932779", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @("BEE ONE HEATER 10138", "https://mrawholesale.com/products/BEE-ONE-HEATER-10138-p598951019", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @('5.5" HEAVY PIPE 10139', "https://mrawholesale.com/products/5-5-HEAVY-PIPE-10139-p598951037", "n/a", "['n/a']", "https://mrawholesale.com/"),
    @('5" HEAVY PIPE WITH SMALL HANDLE 10140', "https://mrawholesale.com/products/5-HEAVY-PIPE-WITH-SMALL-HANDLE-10140-p598950090", "n/a", "['n/a']", "https://mrawholesale.com/")
)

$rowIndex = 9
foreach ($row in $newRows) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
